# Import term 'bead' into the OBI_input mapping sheet (Sheet1).
# A new row is inserted at row 230 (pushing the existing rows 230-232 down
# to 231-233), containing:
#   A230 = http://purl.obolibrary.org/obo/OBI_1000207  (source ontology term IRI)
#   B230 = bead                                        (source ontology term label)
#   C230 = y                                            (Include in View)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above the current row 230, shifting rows down.
$ws.Rows.Item(230).Insert()

# Populate the new row. "bead" is written before the IRI so the new shared
# strings are appended in that order (matches the workbook's shared string table).
$ws.Cells.Item(230, 2).Value = "bead"
$ws.Cells.Item(230, 1).Value = "http://purl.obolibrary.org/obo/OBI_1000207"
$ws.Cells.Item(230, 3).Value = "y"

# Restore the view state (scroll position / selection) reflected in the saved file.
$excel.ActiveWindow.ScrollRow = 223
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("C237").Select()
